$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aggre")

$data = @(
    @("ComoSeEnteroCitaAg", "ComoSeEnteroCitaAg"),
    @("FechaAlta", "High Date"),
    @("IdAgenda", "Agenda ID"),
    @("IdDeteccionMA", "IdDeteccionMA"),
    @("IdDeteccionRE", "IdDeteccionRE"),
    @("IdEstatusCita", "Appointment Status"),
    @("Referencia", "Reference"),
    @("Resultado", "Outcome")
)

$startRow = 332
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

$endRow = $startRow + $data.Length - 1
$range = $ws.Range("A$startRow`:A$endRow")
$fc = $range.FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1
$fc.Font.Color = 393372
$fc.Interior.Color = 13551615
